$d = $word.ActiveDocument
$dash = [char]8211

function Normalize-ParaText($t) {
    # Range.Text can carry a trailing cell-mark (0x07) and/or
    # paragraph-mark (0x0D); strip them so comparisons are exact.
    $t = $t.TrimEnd([char]7)
    $t = $t.TrimEnd([char]13)
    return $t
}

# Returns the 1-based index (in $d.Paragraphs) of the paragraph whose
# normalized text equals $text. When more than one paragraph has that
# text, $nextStartsWith disambiguates by requiring the following
# paragraph's normalized text to start with the given string.
function Get-ParaIndex($d, $text, $nextStartsWith = $null) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $actual = Normalize-ParaText $p.Range.Text
        if ($actual -eq $text) {
            if ($null -eq $nextStartsWith) {
                return $i
            }
            if ($i -lt $d.Paragraphs.Count) {
                $next = $d.Paragraphs.Item($i + 1)
                $nextActual = Normalize-ParaText $next.Range.Text
                if ($nextActual.StartsWith($nextStartsWith)) {
                    return $i
                }
            }
        }
    }
    return -1
}

# ---------------------------------------------------------------------
# Hunk 1: row "1/13", time column -- after the "2:50 - 3:15 PM"
# paragraph (the one immediately followed by "Can, in fact..."),
# add a blank paragraph and a new paragraph "10:30 - 11:15 PM".
# ---------------------------------------------------------------------
$target = "2:50 " + $dash + " 3:15 PM"
$idx = Get-ParaIndex $d $target "Can, in fact"
$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphAfter()
$p2idx = $idx + 1
$p2 = $d.Paragraphs.Item($p2idx)
$p2.Range.InsertParagraphAfter()
$p3idx = $idx + 2
$p3 = $d.Paragraphs.Item($p3idx)
$p3.Range.InsertAfter("10:30 – 11:15 PM")

# ---------------------------------------------------------------------
# Hunk 2: row "1/13", accomplished column -- append a "." to the
# "Looked at Events using ScriptableObjects more closely" paragraph,
# then add a new paragraph with further notes.
# ---------------------------------------------------------------------
$idx = Get-ParaIndex $d "Looked at Events using ScriptableObjects more closely"
$p = $d.Paragraphs.Item($idx)
$p.Range.InsertAfter(".")
$p.Range.InsertParagraphAfter()
$p2idx = $idx + 1
$p2 = $d.Paragraphs.Item($p2idx)
$p2.Range.InsertAfter("Read the Unity Documentation for Serializing and ScriptableObjects, tried filling in more of ScriptableObject_attempt")

# ---------------------------------------------------------------------
# Hunk 3: row "1/14" -- fill in the previously-empty time column and
# accomplished column.
# ---------------------------------------------------------------------
$idx = Get-ParaIndex $d "1/14"
$timeIdx = $idx + 1
$accIdx = $idx + 2
$pTime = $d.Paragraphs.Item($timeIdx)
$pTime.Range.InsertAfter("9:05 – 9:50 AM")

$pAcc = $d.Paragraphs.Item($accIdx)
$rAcc = $pAcc.Range
$rAcc.InsertAfter("Even with the KeyCode method, th")
$rAcc.Collapse(0)
$rAcc.InsertAfter("e movement issue still persists (looking at the Unite 2017 project on v 2020, must be a keyboard problem or something like that). ")
$rAcc.Collapse(0)
$rAcc.InsertAfter("The operator keyword is used after the implicit/explicit keyword, before the return type, with the method having no name; calls the method when the variable with the Object is called without ")
$rAcc.Collapse(0)
$rAcc.InsertAfter("a specific method, returning the return type after the operation stuff is done.")

# ---------------------------------------------------------------------
# Hunk 4: row "1/15", date column -- a `lastRenderedPageBreak` marker
# is recorded just before the "1/" run. This is a pagination-cache
# artifact that Word stamps in during its own layout pass; it carries
# no text/content and is not reachable through the exposed object
# model, so it is intentionally not reproduced here.
# ---------------------------------------------------------------------
